$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values (rows 3-33) to the new figures from the dataset.
$newValues = @{
    3  = 1
    4  = 1
    5  = 12
    6  = 20
    7  = 15
    8  = 16
    9  = 11
    10 = 15
    11 = 21
    12 = 10
    13 = 9
    14 = 0
    15 = 4
    16 = 6
    17 = 12
    18 = 9
    19 = 4
    20 = 8
    21 = 5
    22 = 2
    23 = 7
    24 = 3
    25 = 1
    26 = 1
    27 = 0
    28 = 3
    29 = 1
    30 = 3
    31 = 8
    32 = 1
    33 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("B$row").Value = $newValues[$row]
}

# Update the active selection to match the saved view state.
$ws.Range("E3").Select()
